$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/hour data per latest scrape (GitHub Actions run).
# Values are stored as text in the sheet (inline strings), so force the Text number
# format before assigning to avoid Excel auto-converting "309.89" / "1.01%" / "22" to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.01%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "22"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.31%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "22"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.127"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.96%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "22"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07665"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.01%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "22"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.620"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.58%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "22"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9167"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.26%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "22"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.465"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.48%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "22"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1240"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "21.43%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "22"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1806"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.16%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "22"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09100"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.22%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "22"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04251"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.86%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "22"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.26%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "22"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001248"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.00%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "22"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005635"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.34%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "22"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.358"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.25%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "22"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.283"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.08%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "22"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3313"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.41%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "22"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.911"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.20%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "22"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.71%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "22"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2731"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.02%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "22"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04036"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.73%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "22"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001269"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.24%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "22"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004132"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.93%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "22"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001269"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.69%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "22"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "24.39%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "22"

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "22"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "22"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "22"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "22"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "22"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "22"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "22"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "22"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "22"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "22"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "22"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02432"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.90%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "22"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05256"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.14%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "22"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007830"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.55%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "22"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1309"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.21%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "22"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006800"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.79%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "22"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001841"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.32%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "22"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008190"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.82%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "22"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3341"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.07%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "22"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006857"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.58%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "22"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.36%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "22"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1633"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2,494.76%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "22"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004096"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.16%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "22"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.36%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "22"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.36%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "22"
